$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07631500000000001
$ws.Range("H2").Value = 0.228945
$ws.Range("I2").Value = 0.3245462702943307
$ws.Range("J2").Value = 0.3245462702943307
$ws.Range("M2").Value = 1.949849666666667
$ws.Range("N2").Value = 5.849549000000001
$ws.Range("O2").Value = 0.06676506732104066
$ws.Range("P2").Value = 0.06676506732104066
$ws.Range("Q2").Value = 0.1488027773116667
$ws.Range("R2").Value = 1.339224995805
$ws.Range("S2").Value = 0.02166835358499365
$ws.Range("T2").Value = 0.02166835358499365

# Row 3
$ws.Range("G3").Value = 0.07631500000000001
$ws.Range("H3").Value = 0.228945
$ws.Range("I3").Value = 0.3245462702943307
$ws.Range("J3").Value = 0.3245462702943307
$ws.Range("O3").Value = 0.7967262871802238
$ws.Range("P3").Value = 0.7967262871802239
$ws.Range("Q3").Value = 1.775705305883333
$ws.Range("R3").Value = 15.98134775295
$ws.Range("S3").Value = 0.2585745449497914
$ws.Range("T3").Value = 0.2585745449497915

# Row 4
$ws.Range("G4").Value = 0.07631500000000001
$ws.Range("H4").Value = 0.228945
$ws.Range("I4").Value = 0.3245462702943307
$ws.Range("J4").Value = 0.3245462702943307
$ws.Range("O4").Value = 0.1365086454987356
$ws.Range("P4").Value = 0.1365086454987356
$ws.Range("Q4").Value = 0.3042439166516667
$ws.Range("R4").Value = 2.738195249865
$ws.Range("S4").Value = 0.04430337175954562
$ws.Range("T4").Value = 0.04430337175954562

# Row 5
$ws.Range("G5").Value = 0.1588286666666667
$ws.Range("I5").Value = 0.6754537297056692
$ws.Range("J5").Value = 0.6754537297056693
$ws.Range("M5").Value = 1.949849666666667
$ws.Range("N5").Value = 5.849549000000001
$ws.Range("O5").Value = 0.06676506732104066
$ws.Range("P5").Value = 0.06676506732104066
$ws.Range("Q5").Value = 0.3096920227571112
$ws.Range("R5").Value = 2.787228204814
$ws.Range("S5").Value = 0.045096713736047
$ws.Range("T5").Value = 0.04509671373604701

# Row 6
$ws.Range("G6").Value = 0.1588286666666667
$ws.Range("I6").Value = 0.6754537297056692
$ws.Range("J6").Value = 0.6754537297056693
$ws.Range("O6").Value = 0.7967262871802238
$ws.Range("P6").Value = 0.7967262871802239
$ws.Range("S6").Value = 0.5381517422304323
$ws.Range("T6").Value = 0.5381517422304324

# Row 7
$ws.Range("G7").Value = 0.1588286666666667
$ws.Range("I7").Value = 0.6754537297056692
$ws.Range("J7").Value = 0.6754537297056693
$ws.Range("O7").Value = 0.1365086454987356
$ws.Range("P7").Value = 0.1365086454987356
$ws.Range("Q7").Value = 0.6331999688557778
$ws.Range("S7").Value = 0.09220527373918999
$ws.Range("T7").Value = 0.09220527373919001
